$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text so numeric-looking strings
# (e.g. "220.54", "26.316.78") are not auto-coerced into numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.316.78'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '1.666.94'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').Value = '220.54'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').Value = '0.5307'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('D9').Value = '0.06369'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('D11').Value = '0.07842'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '4.515'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = '1.675.48'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').Value = '1.895.41'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').Value = '0.5595'
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('D16').Value = '0.0₅8160'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '26.321.90'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').Value = '1.009'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').Value = '4.715'
$ws.Range('D21').Value = '197.14'
$ws.Range('E21').Value = '  +3.17%  '
$ws.Range('D22').Value = '10.26'
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('D23').Value = '6.046'
$ws.Range('E23').Value = '  +0.56%  '
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').Value = '145.96'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').Value = '0.1221'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').Value = '7.239'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('D28').Value = '16.15'
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('D29').Value = '1.506'
$ws.Range('E29').Value = '  +2.12%  '
$ws.Range('D30').Value = '0.05889'
$ws.Range('E30').Value = '  +2.45%  '
$ws.Range('D31').Value = '1.285'
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').Value = '3.539'
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('D33').Value = '3.334'
$ws.Range('E33').Value = '  +2.18%  '
$ws.Range('D34').Value = '1.601'
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').Value = '2.828'
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D36').Value = '0.9598'
$ws.Range('E36').Value = '  +1.15%  '
$ws.Range('D37').Value = '2.434'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('E39').Value = '  +0.72%  '
$ws.Range('D40').Value = '5.944'
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('D41').Value = '1.077.47'
$ws.Range('E41').Value = '  +3.52%  '
$ws.Range('D42').Value = '0.8578'
$ws.Range('E42').Value = '  +0.80%  '
$ws.Range('D44').Value = '102.74'
$ws.Range('E44').Value = '  -0.97%  '
$ws.Range('D45').Value = '1.805.86'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('E46').Value = '  +2.95%  '
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('D50').Value = '8.014'
$ws.Range('E50').Value = '  +2.20%  '
$ws.Range('D51').Value = '0.05150'
$ws.Range('E51').Value = '  -0.05%  '

# Restore default (General) formatting now that the text values are set,
# so no stray number-format style is left behind on these cells.
$dataRange.ClearFormats()

Write-Host "Applied 84 cell updates"
